$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29; this shifts existing rows 29..129 down to 30..130
# and automatically extends the used range / dimension to row 130.
$ws.Rows(29).Insert()

# Populate the newly inserted row 29 with the same record as the (now shifted down)
# row 30, except for the date (column D) and volume (column J), which are the
# genuinely new values for this weekly entry.
$ws.Cells.Item(29, 1).Value = 7
$ws.Cells.Item(29, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(29, 3).Value = "Ñuble"
$ws.Cells.Item(29, 4).Value = 44487
$ws.Cells.Item(29, 5).Value = 16
$ws.Cells.Item(29, 6).Value = 100112017
$ws.Cells.Item(29, 7).Value = "Apio"
$ws.Cells.Item(29, 8).Value = "Americana (o)"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 180
$ws.Cells.Item(29, 11).Value = 8000
$ws.Cells.Item(29, 12).Value = 8500
$ws.Cells.Item(29, 13).Value = 8250
$ws.Cells.Item(29, 14).Value = "$/docena de matas"
$ws.Cells.Item(29, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(29, 16).Value = 1375
$ws.Cells.Item(29, 17).Value = 6
$ws.Cells.Item(29, 18).Value = "Hortaliza"
